$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

# New values for B2:G11 (row numbers correspond to sheet rows)
$data = @{
    2  = @(-0.3684425897240435, 1.990300204684291, 12.34821194374062, 3.514002268602088, 3.576871232156456, 22)
    3  = @(-0.4406390919636032, 1.050095607268688, 4.94590582631831, 2.223939258684533, 2.233680856299617, 21)
    4  = @(-0.2266117036152265, 0.7762331526076556, 2.687494923465916, 1.639358082746389, 1.665799020938133, 20)
    5  = @(0.02747715458826221, 0.7115633747814063, 1.229190786505692, 1.10868876899953, 1.138719555200954, 19)
    6  = @(-0.06432496402825524, 0.6672897034160007, 1.122755052931566, 1.059601365104616, 1.088309869768877, 18)
    7  = @(0.1299503236878231, 0.5419393736802232, 0.7287712587988439, 0.8536810052934549, 0.8696993347893518, 17)
    8  = @(0.1289229671324676, 0.5250093336297077, 0.4739960288771982, 0.6884736951236395, 0.6984744499630291, 16)
    9  = @(0.196164091590951, 0.5034467950485412, 0.4366068896125055, 0.6607623548693626, 0.6531187860533945, 15)
    10 = @(0.2421836099529805, 0.4042897487511507, 0.229010732524476, 0.4785506582635489, 0.4283249702957092, 14)
    11 = @(0.2725419345940574, 0.3839238355923389, 0.1922556345917136, 0.4384696507076785, 0.3575023158704416, 13)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$rowNum").Value = $vals[$i]
    }
}
